$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 44
$ws.Range("I2").Value = 84
$ws.Range("J2").Value = 366
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 114
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = 75
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 2
$ws.Range("T2").Value = 89
$ws.Range("V2").Value = 610
$ws.Range("X2").Value = 685
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 6
$ws.Range("AA2").Value = 4
